$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ER (Etienne) - "Developpement Web 3" on Lundi: 13:15-16:05 -> 12:15-15:05
$ws.Range("B14").Value = "'12:15"
$ws.Range("C14").Value = "'15:05"

# ER (Etienne) - "Developpement Web 3" on Jeudi: 10:15-12:05 -> 13:15-15:05
$ws.Range("B21").Value = "'13:15"
$ws.Range("C21").Value = "'15:05"

# AO (Alexandre) - "Projet integrateur 1" on Lundi: 10:15-13:05 -> 8:15-11:05
$ws.Range("B28").Value = "'8:15"
$ws.Range("C28").Value = "'11:05"

# AO (Alexandre) - "Projet integrateur 1" moved from Mercredi 12:15-14:05 to Vendredi 15:15-17:05
$ws.Range("A29").Value = "Vendredi"
$ws.Range("B29").Value = "'15:15"
$ws.Range("C29").Value = "'17:05"

# ER (Etienne) - "Piratage ethique" moved from Vendredi 10:15-12:05 (C211) to Vendredi 15:15-17:05 (E209)
$ws.Range("B32").Value = "'15:15"
$ws.Range("C32").Value = "'17:05"
$ws.Range("G32").Value = "E209"

# Re-sort the table by person (column E) then by day (column A), ascending - matches Table1's sortState
$tbl = $ws.ListObjects.Item(1)
$sort = $tbl.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("E2:E39"))
$sort.SortFields.Add($ws.Range("A2:A39"))
$sort.Header = 1
$sort.Apply()

# Leave selection on G10 (matches the final saved view state)
$ws.Range("G10").Select()
